# Add carjacking data for 2021-09-29 (one additional day of data)
# across the "through September <day>" columns for each year, plus a few
# other newly-observed incidents scattered across the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / update the "through" date in the title row.
$ws.Name = "Through 2021-09-29"
$ws.Range("B1").Value = "September 2021 (through September 29)"

# Updated counts: incremented totals for the day added.
$ws.Range("B2").Value = 17    # Garfield Park, Sep 2021
$ws.Range("AL2").Value = 7    # Garfield Park, Sep 2017
$ws.Range("BD2").Value = 4    # Garfield Park, Sep 2015

$ws.Range("K3").Value = 9     # North Lawndale, Sep 2020
$ws.Range("T3").Value = 6     # North Lawndale, Sep 2019
$ws.Range("AC3").Value = 3    # North Lawndale, Sep 2018

$ws.Range("AL4").Value = 2    # Humboldt Park, Sep 2017

$ws.Range("AU5").Value = 3    # Austin, Sep 2016

$ws.Range("K7").Value = 3     # Auburn Gresham, Sep 2020
$ws.Range("AU7").Value = 1    # Auburn Gresham, Sep 2016 (new)

$ws.Range("T13").Value = 4    # Chatham, Sep 2019

$ws.Range("B20").Value = 2    # Englewood, Sep 2021

$ws.Range("B27").Value = 4    # Avalon Park, Sep 2021

$ws.Range("B31").Value = 2    # West Loop, Sep 2021

$ws.Range("K33").Value = 3    # Lake View, Sep 2020
$ws.Range("K34").Value = 2    # Hyde Park, Sep 2020 (new)

$ws.Range("BD39").Value = 1   # Douglas, Sep 2015 (new)

$ws.Range("AC53").Value = 1   # Irving Park, Sep 2018 (new)

$ws.Range("B94").Value = 2    # Rush & Division, Sep 2021

$ws.Range("B99").Value = 1    # West Ridge, Sep 2021 (new)
